$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after row 8 to hold the data that used to live in row 8
$ws.Rows.Item(9).Insert()

# New row 9 gets the values that were previously in row 8 (before this edit)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44838
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112039
$ws.Range("G9").Value = "Ciboulette"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 1000
$ws.Range("N9").Value = "$/docena de atados"
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 333
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = "Hortaliza"

# Row 8 updates (quality moves from Segunda to Primera with row7's old numbers)
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 1200
$ws.Range("L8").Value = 1300
$ws.Range("M8").Value = 1250
$ws.Range("P8").Value = 417

# Row 7 updates (new date + new numbers)
$ws.Range("D7").Value = 45134
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("P7").Value = 833
